# candidateGenes.xlsx update:
# - "Flag"/"Azoo" pathology labels are renamed to their proper cohort
#   acronyms (MMAF / NOA), and the MEI4 candidate gene row is repointed
#   to DNAH1 with a NOA->MMAF pathology re-classification + updated Level.
# - Selection cursor is left on A6 (first empty row below the table),
#   matching where the author's cursor ended up after editing the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: CFAP44 / Flag -> MMAF / Level 5 (unchanged)
$ws.Range("B2").Value = "MMAF"

# Row 3: MEI1 / Azoo -> NOA / Level 4 -> 3
$ws.Range("B3").Value = "NOA"
$ws.Range("C3").Value = 3

# Row 4: MEI4/Azoo -> DNAH1/MMAF / Level 3 -> 4
$ws.Range("A4").Value = "DNAH1"
$ws.Range("B4").Value = "MMAF"
$ws.Range("C4").Value = 4

# Row 5: SPINK2 / Azoo -> NOA / Level 5 (unchanged)
$ws.Range("B5").Value = "NOA"

# Move the active selection to A6 (below the last data row)
$ws.Range("A6").Select()
